$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.0630466474718978
$ws.Range("J2").Value = 0.06304664747189781
$ws.Range("M2").Value = 3.770298333333333
$ws.Range("N2").Value = 11.310895
$ws.Range("O2").Value = 0.06100259562224731
$ws.Range("P2").Value = 0.06125631726190612
$ws.Range("Q2").Value = 1.075769169321111
$ws.Range("R2").Value = 9.68192252389
$ws.Range("S2").Value = 0.003846009141066562
$ws.Range("T2").Value = 0.003862005439838124

# Row 3
$ws.Range("I3").Value = 0.0630466474718978
$ws.Range("J3").Value = 0.06304664747189781
$ws.Range("O3").Value = 0.06469423882843597
$ws.Range("P3").Value = 0.06496331472897099
$ws.Range("S3").Value = 0.004078754868879165
$ws.Range("T3").Value = 0.004095719202323381

# Row 4
$ws.Range("I4").Value = 0.0630466474718978
$ws.Range("J4").Value = 0.06304664747189781
$ws.Range("M4").Value = 31.40746233333333
$ws.Range("N4").Value = 94.222387
$ws.Range("O4").Value = 0.5081658147055464
$ws.Range("P4").Value = 0.5102793749960634
$ws.Range("Q4").Value = 8.961407474337109
$ws.Range("R4").Value = 80.652667269034
$ws.Range("S4").Value = 0.03203815097701032
$ws.Range("T4").Value = 0.03217140386755715

# Row 5
$ws.Range("I5").Value = 0.0630466474718978
$ws.Range("J5").Value = 0.06304664747189781
$ws.Range("M5").Value = 0.7679895
$ws.Range("N5").Value = 1.535979
$ws.Range("O5").Value = 0.01242590075603175
$ws.Range("P5").Value = 0.008318388326620067
$ws.Range("Q5").Value = 0.219128396063
$ws.Range("R5").Value = 1.314770376378
$ws.Range("S5").Value = 0.0007834113844863218
$ws.Range("T5").Value = 0.0005244464963627654

# Row 6
$ws.Range("I6").Value = 0.0630466474718978
$ws.Range("J6").Value = 0.06304664747189781
$ws.Range("M6").Value = 21.861327
$ws.Range("N6").Value = 65.58398100000001
$ws.Range("O6").Value = 0.3537114500877385
$ws.Range("P6").Value = 0.3551826046864394
$ws.Range("Q6").Value = 6.237634136038
$ws.Range("R6").Value = 56.13870722434201
$ws.Range("S6").Value = 0.02230032110045542
$ws.Range("T6").Value = 0.02239307246581638

# Row 7
$ws.Range("G7").Value = 3.253975333333333
$ws.Range("H7").Value = 9.761925999999999
$ws.Range("I7").Value = 0.7190065996349845
$ws.Range("J7").Value = 0.7190065996349846
$ws.Range("M7").Value = 3.770298333333333
$ws.Range("N7").Value = 11.310895
$ws.Range("O7").Value = 0.06100259562224731
$ws.Range("P7").Value = 0.06125631726190612
$ws.Range("Q7").Value = 12.26845777597444
$ws.Range("R7").Value = 110.41611998377
$ws.Range("S7").Value = 0.04386126884726003
$ws.Range("T7").Value = 0.04404369638064493

# Row 8
$ws.Range("G8").Value = 3.253975333333333
$ws.Range("H8").Value = 9.761925999999999
$ws.Range("I8").Value = 0.7190065996349845
$ws.Range("J8").Value = 0.7190065996349846
$ws.Range("O8").Value = 0.06469423882843597
$ws.Range("P8").Value = 0.06496331472897099
$ws.Range("Q8").Value = 13.01089780392911
$ws.Range("R8").Value = 117.098080235362
$ws.Range("S8").Value = 0.04651558467600733
$ws.Range("T8").Value = 0.04670905202429474

# Row 9
$ws.Range("G9").Value = 3.253975333333333
$ws.Range("H9").Value = 9.761925999999999
$ws.Range("I9").Value = 0.7190065996349845
$ws.Range("J9").Value = 0.7190065996349846
$ws.Range("M9").Value = 31.40746233333333
$ws.Range("N9").Value = 94.222387
$ws.Range("O9").Value = 0.5081658147055464
$ws.Range("P9").Value = 0.5102793749960634
$ws.Range("Q9").Value = 102.1991077152624
$ws.Range("R9").Value = 919.7919694373619
$ws.Range("S9").Value = 0.3653745744821765
$ws.Range("T9").Value = 0.3668942382797847

# Row 10
$ws.Range("G10").Value = 3.253975333333333
$ws.Range("H10").Value = 9.761925999999999
$ws.Range("I10").Value = 0.7190065996349845
$ws.Range("J10").Value = 0.7190065996349846
$ws.Range("M10").Value = 0.7679895
$ws.Range("N10").Value = 1.535979
$ws.Range("O10").Value = 0.01242590075603175
$ws.Range("P10").Value = 0.008318388326620067
$ws.Range("Q10").Value = 2.499018889259
$ws.Range("R10").Value = 14.994113335554
$ws.Range("S10").Value = 0.008934304649996169
$ws.Range("T10").Value = 0.005980976105166444

# Row 11
$ws.Range("G11").Value = 3.253975333333333
$ws.Range("H11").Value = 9.761925999999999
$ws.Range("I11").Value = 0.7190065996349845
$ws.Range("J11").Value = 0.7190065996349846
$ws.Range("M11").Value = 21.861327
$ws.Range("N11").Value = 65.58398100000001
$ws.Range("O11").Value = 0.3537114500877385
$ws.Range("P11").Value = 0.3551826046864394
$ws.Range("Q11").Value = 71.136218811934
$ws.Range("R11").Value = 640.225969307406
$ws.Range("S11").Value = 0.2543208669795444
$ws.Range("T11").Value = 0.2553786368450937

# Row 12
$ws.Range("G12").Value = 0.9863516666666667
$ws.Range("H12").Value = 2.959055
$ws.Range("I12").Value = 0.2179467528931175
$ws.Range("J12").Value = 0.2179467528931176
$ws.Range("M12").Value = 3.770298333333333
$ws.Range("N12").Value = 11.310895
$ws.Range("O12").Value = 0.06100259562224731
$ws.Range("P12").Value = 0.06125631726190612
$ws.Range("Q12").Value = 3.718840044913889
$ws.Range("R12").Value = 33.46956040422501
$ws.Range("S12").Value = 0.01329531763392071
$ws.Range("T12").Value = 0.01335061544142307

# Row 13
$ws.Range("G13").Value = 0.9863516666666667
$ws.Range("H13").Value = 2.959055
$ws.Range("I13").Value = 0.2179467528931175
$ws.Range("J13").Value = 0.2179467528931176
$ws.Range("O13").Value = 0.06469423882843597
$ws.Range("P13").Value = 0.06496331472897099
$ws.Range("Q13").Value = 3.943889986587223
$ws.Range("R13").Value = 35.49500987928501
$ws.Range("S13").Value = 0.01409989928354947
$ws.Range("T13").Value = 0.01415854350235287

# Row 14
$ws.Range("G14").Value = 0.9863516666666667
$ws.Range("H14").Value = 2.959055
$ws.Range("I14").Value = 0.2179467528931175
$ws.Range("J14").Value = 0.2179467528931176
$ws.Range("M14").Value = 31.40746233333333
$ws.Range("N14").Value = 94.222387
$ws.Range("O14").Value = 0.5081658147055464
$ws.Range("P14").Value = 0.5102793749960634
$ws.Range("Q14").Value = 30.97880281825389
$ws.Range("R14").Value = 278.809225364285
$ws.Range("S14").Value = 0.1107530892463595
$ws.Range("T14").Value = 0.1112137328487215

# Row 15
$ws.Range("G15").Value = 0.9863516666666667
$ws.Range("H15").Value = 2.959055
$ws.Range("I15").Value = 0.2179467528931175
$ws.Range("J15").Value = 0.2179467528931176
$ws.Range("M15").Value = 0.7679895
$ws.Range("N15").Value = 1.535979
$ws.Range("O15").Value = 0.01242590075603175
$ws.Range("P15").Value = 0.008318388326620067
$ws.Range("Q15").Value = 0.7575077233075
$ws.Range("R15").Value = 4.545046339845
$ws.Range("S15").Value = 0.002708184721549253
$ws.Range("T15").Value = 0.001812965725090858

# Row 16
$ws.Range("G16").Value = 0.9863516666666667
$ws.Range("H16").Value = 2.959055
$ws.Range("I16").Value = 0.2179467528931175
$ws.Range("J16").Value = 0.2179467528931176
$ws.Range("M16").Value = 21.861327
$ws.Range("N16").Value = 65.58398100000001
$ws.Range("O16").Value = 0.3537114500877385
$ws.Range("P16").Value = 0.3551826046864394
$ws.Range("Q16").Value = 21.562956321995
$ws.Range("R16").Value = 194.066606897955
$ws.Range("S16").Value = 0.07709026200773862
$ws.Range("T16").Value = 0.07741089537552927
